$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (E1) by deleting the cell and shifting
# the remaining header cells to the left, preserving formatting.
$ws.Range("E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
